$d = $word.ActiveDocument

# Remove the two paragraphs that follow "As reuniões vão ser presenciais...":
#   "Contato Barbearia Brothers:"                      (paragraph w/ hyperlink)
#   "Telefone (WhatsApp): (11) 96727-9523"
# Locate them by their distinctive text and delete the full paragraph
# ranges (including paragraph marks) so the following empty paragraph
# collapses back next to the "As reuniões..." paragraph, matching the diff.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "Contato*") {
        $startPara = $p
    }
    if ($t -like "Telefone (WhatsApp)*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
